$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 92179
$ws.Range("B3").Value = 89193
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 57881
$ws.Range("B6").Value = 57881
$ws.Range("B7").Value = 91808
$ws.Range("A8").Value = 131046711
$ws.Range("B8").Value = 83223
$ws.Range("E8").Value = 6440
$ws.Range("F8").Value = "Vitgrynig nållav"
$ws.Range("G8").Value = "Chaenotheca subroscida"
$ws.Range("H8").Value = "(Eitner) Zahlbr."
$ws.Range("Q8").Value = 402363
$ws.Range("R8").Value = 6818428
$ws.Range("Z8").Value = "16:09"
$ws.Range("AB8").Value = "16:09"
$ws.Range("A9").Value = 131046843
$ws.Range("B9").Value = 79243
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 402432
$ws.Range("R9").Value = 6818480
$ws.Range("Z9").Value = "16:20"
$ws.Range("AB9").Value = "16:20"
$ws.Range("A10").Value = 131046844
$ws.Range("B10").Value = 79243
$ws.Range("Q10").Value = 402484
$ws.Range("R10").Value = 6818538
$ws.Range("Z10").Value = "16:23"
$ws.Range("AB10").Value = "16:23"
$ws.Range("A11").Value = 131046763
$ws.Range("B11").Value = 92267
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 1209
$ws.Range("F11").Value = "Rynkskinn"
$ws.Range("G11").Value = "Hermanssonia centrifuga"
$ws.Range("H11").Value = "(P. Karst.) Zmitr."
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 402378
$ws.Range("R11").Value = 6818392
$ws.Range("Z11").Value = "17:02"
$ws.Range("AB11").Value = "17:02"
$ws.Range("AC11").ClearContents()
$ws.Range("A12").Value = 131046788
$ws.Range("B12").Value = 57884
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 402473
$ws.Range("R12").Value = 6818425
$ws.Range("Z12").Value = "16:47"
$ws.Range("AB12").Value = "16:47"
$ws.Range("AC12").Value = "Färska ringhack (gran)"
$ws.Range("B13").Value = 57884
$ws.Range("A14").Value = 131046806
$ws.Range("B14").Value = 83206
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 6439
$ws.Range("F14").Value = "Gulnål"
$ws.Range("G14").Value = "Chaenotheca brachypoda"
$ws.Range("H14").Value = "(Ach.) Tibell"
$ws.Range("Q14").Value = 402340
$ws.Range("R14").Value = 6818363
$ws.Range("Z14").Value = "17:05"
$ws.Range("AB14").Value = "17:05"
$ws.Range("A15").Value = 131046811
$ws.Range("B15").Value = 91828
$ws.Range("Q15").Value = 402450
$ws.Range("R15").Value = 6818298
$ws.Range("Z15").Value = "16:54"
$ws.Range("AB15").Value = "16:54"
$ws.Range("A16").Value = 131046808
$ws.Range("B16").Value = 91828
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = "Granticka"
$ws.Range("G16").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H16").Value = ""
$ws.Range("Q16").Value = 402323
$ws.Range("R16").Value = 6818416
$ws.Range("Z16").Value = "16:06"
$ws.Range("AB16").Value = "16:06"
$ws.Range("B17").Value = 57073
$ws.Range("B18").Value = 91828
$ws.Range("B19").Value = 57881
$ws.Range("B20").Value = 79243
$ws.Range("B21").Value = 57881
$ws.Range("B22").Value = 91771
$ws.Range("A23").Value = 131046845
$ws.Range("B23").Value = 79243
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("M23").ClearContents()
$ws.Range("Q23").Value = 402575
$ws.Range("R23").Value = 6818545
$ws.Range("Z23").Value = "16:34"
$ws.Range("AB23").Value = "16:34"
$ws.Range("AC23").ClearContents()
$ws.Range("AE23").Value = $False
$ws.Range("A24").Value = 131047016
$ws.Range("B24").Value = 57884
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("M24").Value = "färska spår"
$ws.Range("Q24").Value = 402474
$ws.Range("R24").Value = 6818507
$ws.Range("Z24").Value = "16:22"
$ws.Range("AB24").Value = "16:22"
$ws.Range("AC24").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE24").Value = $True
$ws.Range("B25").Value = 79243
